# RippleTemplate_MultiSrcConc.xlsx update
#  - add a new "Assay" worksheet (settings/value table) at the end of the tab
#    order
#  - make "Patterns" the selected/active tab (was "Compounds")
#  - tidy up a couple of selection anchors along the way

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. New "Assay" sheet, appended after the last existing sheet ("Barcodes")
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$assay = $wb.Worksheets.Add($null, $lastSheet)
$assay.Name = "Assay"

$assay.Range("A1").Value = "Setting"
$assay.Range("B1").Value = "Value"

$assay.Range("A2").Value = "DMSO Tolerance"
$assay.Range("B2").Value = 0.005

$assay.Range("A3").Value = "Well Volume (µL)"
$assay.Range("B3").Value = 25

$assay.Range("A4").Value = "Backfill (µL)"
$assay.Range("B4").Value = 10

$assay.Range("A5").Value = "Allowed Error"
$assay.Range("B5").Value = 0.1

$assay.Range("A6").Value = "Destination Replicates"
$assay.Range("B6").Value = 1

$assay.Range("A7").Value = "Use Intermediate Plates"
$assay.Range("B7").Value = 1

$assay.Range("A8").Value = "DMSO Normalization"
$assay.Range("B8").Value = 1

[void]$assay.Range("A1:B8").Select()

# ---------------------------------------------------------------------------
# 2. Tab / selection bookkeeping:
#      "Patterns"  becomes the active/selected tab (previously "Compounds"),
#      with the cursor resting on K21 (was E15).
#      "Compounds" loses its tabSelected flag but keeps its prior selection.
# ---------------------------------------------------------------------------
$patterns = $wb.Worksheets.Item("Patterns")
[void]$patterns.Activate()
[void]$patterns.Range("K21").Select()
